$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- last_edited_time (column D) updates -----------------------------
# Rows 2-10 previously shared the "2024-08-26T17:26:00.000Z" timestamp;
# that timestamp is corrected to "2024-08-27T12:18:00.000Z".
$ws.Range("D2:D10").Value = "2024-08-27T12:18:00.000Z"

# Rows 11-22 previously shared (or are re-pointed to) the
# "2024-08-26T17:27:00.000Z" timestamp; corrected to
# "2024-08-27T12:15:00.000Z".
$ws.Range("D11:D22").Value = "2024-08-27T12:15:00.000Z"

# --- chiết khấu (discount) number corrections -------------------------
# Row 10
$ws.Range("S10").Value = 21
$ws.Range("AF10").Value = 20
$ws.Range("AM10").Value = 21

# Row 13
$ws.Range("S13").Value = 4
$ws.Range("AF13").Value = 34
$ws.Range("AM13").Value = 34
$ws.Range("AX13").Value = 5

# Row 16
$ws.Range("S16").Value = 24.5
$ws.Range("AF16").Value = 23
$ws.Range("AM16").Value = 24.5

# Row 18
$ws.Range("AF18").Value = 15
$ws.Range("AM18").Value = 17.5
$ws.Range("AX18").Value = 1

# Row 19
$ws.Range("S19").Value = 23
$ws.Range("AF19").Value = 22
$ws.Range("AM19").Value = 23

# Row 21
$ws.Range("S21").Value = 23
$ws.Range("AF21").Value = 23
$ws.Range("AM21").Value = 23

# Row 22
$ws.Range("S22").Value = 25
$ws.Range("AF22").Value = 24
$ws.Range("AM22").Value = 25
